$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Before: A1=No, B1=Kode, C1=Program, D1=Tahun Perubahan
# After:  A1=No, B1=Kode Urusan, C1=Kode Program, D1=Program, E1=Tahun Perubahan
#
# Insert a new column at C (pushes old C "Program" -> D, old D "Tahun Perubahan" -> E)
$ws.Columns("C").Insert()

# Rename the existing "Kode" header (column B) to "Kode Urusan"
$ws.Range("B1").Value = "Kode Urusan"

# Fill in the newly inserted column C header
$ws.Range("C1").Value = "Kode Program"

# Match column widths for the two new "Kode ..." columns (best-fit sizing)
$ws.Columns("B").ColumnWidth = 10.5
$ws.Columns("C").ColumnWidth = 11.666666666666666

# Update the active selection to D5 (matches the saved view state)
$null = $ws.Range("D5").Select()
